# Actualización 11 de Mayo - Mañana
$wb = $excel.ActiveWorkbook

# --- Hoja "Estadisticos 2P": actualizar estadísticos de los grupos 6ARHV ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

# Fila 3 (6ARHV / DETERMINA LA NÓMINA...)
$ws2.Range("D3").Value = 6
$ws2.Range("E3").Value = 6
$ws2.Range("F3").Value = 22
$ws2.Range("G3").Value = 78.57

# Fila 4 (6ARHV / DETERMINA REMUNERACIONES...)
$ws2.Range("D4").Value = 8
$ws2.Range("E4").Value = 9
$ws2.Range("F4").Value = 19
$ws2.Range("G4").Value = 67.86
$ws2.Range("H4").Value = 8.4

# --- Hoja "Rescatables": reacomodo de materias reprobadas y nuevo alumno ---
$ws4 = $wb.Worksheets.Item("Rescatables")

# Fila 2 y 3: intercambiar la materia reprobada (mismo alumno CIRUELO MANCILLA)
$ws4.Range("E2").Value = "DETERMINA LA NÓMINA DEL PERSONAL DE LA ORGANIZACIÓN TOMANDO EN CUENTA LA NORMATIVIDAD LABORAL"
$ws4.Range("E3").Value = "DETERMINA REMUNERACIONES DEL PERSONAL EN SITUACIONES EXTRAORDINARIAS"

# Fila 4 y 5: intercambiar la materia reprobada (mismo alumno VERA PAZOS)
$ws4.Range("E4").Value = "DETERMINA REMUNERACIONES DEL PERSONAL EN SITUACIONES EXTRAORDINARIAS"
$ws4.Range("E5").Value = "DETERMINA LA NÓMINA DEL PERSONAL DE LA ORGANIZACIÓN TOMANDO EN CUENTA LA NORMATIVIDAD LABORAL"

# Fila 7: se sustituye el alumno (ahora ZEPEDA MORALES NATHAEL) y reprueba 2 materias
$ws4.Range("A7").Value = 18330051920110
$ws4.Range("B7").Value = "ZEPEDA"
$ws4.Range("C7").Value = "MORALES"
$ws4.Range("D7").Value = "NATHAEL"
$ws4.Range("G7").Value = 2

# Filas 8 y 9: ya no aplican (alumnos dados de baja de rescate)
$ws4.Rows.Item(9).Delete()
$ws4.Rows.Item(8).Delete()
